$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while forcing text type (avoids
# Excel auto-converting numeric-looking strings like "54.04" into real
# numbers), then clear the resulting quote-prefix style so the cells
# style index is left untouched (matches original unstyled data cells).
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.849.63"
$ws.Range("E2").Value = "  -0.04%  "
Set-TextValue $ws.Range("D3") "3.520.00"
$ws.Range("E3").Value = "  -0.47%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue $ws.Range("D5") "601.27"
$ws.Range("E5").Value = "  -1.67%  "
Set-TextValue $ws.Range("D6") "195.42"
$ws.Range("E6").Value = "  +5.68%  "
Set-TextValue $ws.Range("D7") "0.624"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E9").Value = "  -2.21%  "
Set-TextValue $ws.Range("D10") "0.653"
$ws.Range("E10").Value = "  +1.34%  "
Set-TextValue $ws.Range("D11") "54.04"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("E12").Value = "  -2.58%  "
Set-TextValue $ws.Range("D13") "9.52"
$ws.Range("E13").Value = "  +1.01%  "
Set-TextValue $ws.Range("D14") "4.074.14"
Set-TextValue $ws.Range("D15") "603.09"
$ws.Range("E15").Value = "  -2.99%  "
Set-TextValue $ws.Range("D16") "70.038.87"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("E17").Value = "  +1.59%  "
Set-TextValue $ws.Range("D18") "12.62"
$ws.Range("E18").Value = "  +0.00%  "
Set-TextValue $ws.Range("D19") "3.518.42"
$ws.Range("E19").Value = "  -0.65%  "
Set-TextValue $ws.Range("D21") "0.995"
$ws.Range("E21").Value = "  +0.45%  "
Set-TextValue $ws.Range("D22") "18.28"
$ws.Range("E22").Value = "  +4.20%  "
$ws.Range("E23").Value = "  +6.45%  "
Set-TextValue $ws.Range("D24") "103.80"
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("E25").Value = "  -2.45%  "
$ws.Range("E26").Value = "  +2.71%  "
Set-TextValue $ws.Range("D27") "10.92"
$ws.Range("E27").Value = "  -0.29%  "
Set-TextValue $ws.Range("D28") "9.67"
$ws.Range("E28").Value = "  +1.23%  "
Set-TextValue $ws.Range("D29") "33.55"
$ws.Range("E29").Value = "  +2.36%  "
Set-TextValue $ws.Range("D30") "4.49"
$ws.Range("E30").Value = "  +24.58%  "
Set-TextValue $ws.Range("D31") "7.10"
$ws.Range("E31").Value = "  +1.31%  "
Set-TextValue $ws.Range("D32") "12.70"
$ws.Range("E32").Value = "  +3.97%  "
$ws.Range("E33").Value = "  +1.84%  "
Set-TextValue $ws.Range("D34") "63.23"
$ws.Range("E34").Value = "  -0.33%  "
Set-TextValue $ws.Range("D35") "3.765.66"
$ws.Range("E35").Value = "  +6.44%  "
Set-TextValue $ws.Range("D36") "0.0₃0823"
$ws.Range("E36").Value = "  +5.69%  "
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("E40").Value = "  +1.26%  "
Set-TextValue $ws.Range("D41") "36.77"
$ws.Range("E41").Value = "  -0.79%  "
Set-TextValue $ws.Range("D42") "488.08"
$ws.Range("E42").Value = "  -8.53%  "
$ws.Range("E43").Value = "  -0.13%  "
Set-TextValue $ws.Range("D44") "0.0456"
$ws.Range("E44").Value = "  -0.05%  "
Set-TextValue $ws.Range("D45") "0.140"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("E49").Value = "  -5.47%  "

# Row swaps: the two rows trade coin identity (name/link) and carry new
# price/volume data as given by the source diff (not a straight value swap).
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D37") "1.00"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D38") "3.08"
$ws.Range("E38").Value = "  -5.13%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D46") "2.83"
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D47") "3.32"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D50") "1.33"
$ws.Range("E50").Value = "  +13.53%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws.Range("D51") "0.000243"
$ws.Range("E51").Value = "  +1.02%  "
